$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete "ECs sending cluster" rows (old rows 8-10)
$ws.Rows.Item(10).EntireRow.Delete()
$ws.Rows.Item(9).EntireRow.Delete()
$ws.Rows.Item(8).EntireRow.Delete()

# Update remaining data rows (2-7) with the refreshed TPM-derived values
# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Lama1"
$ws.Range("C2").Value = "Itgb8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3395593333333333
$ws.Range("H2").Value = 1.018678
$ws.Range("I2").Value = 0.6166145092460882
$ws.Range("J2").Value = 0.6166145092460882
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1591403333333333
$ws.Range("N2").Value = 0.477421
$ws.Range("O2").Value = 0.01260326986877891
$ws.Range("P2").Value = 0.01260326986877891
$ws.Range("Q2").Value = 0.05403758549311111
$ws.Range("R2").Value = 0.486338269438
$ws.Range("S2").Value = 0.007771359065033117
$ws.Range("T2").Value = 0.007771359065033118

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Lama1"
$ws.Range("C3").Value = "Itgb8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3395593333333333
$ws.Range("H3").Value = 1.018678
$ws.Range("I3").Value = 0.6166145092460882
$ws.Range("J3").Value = 0.6166145092460882
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.794584666666667
$ws.Range("N3").Value = 11.383754
$ws.Range("O3").Value = 0.3005157372251983
$ws.Range("P3").Value = 0.3005157372251983
$ws.Range("Q3").Value = 1.288486639690222
$ws.Range("R3").Value = 11.596379757212
$ws.Range("S3").Value = 0.185302363829842
$ws.Range("T3").Value = 0.1853023638298421

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Lama1"
$ws.Range("C4").Value = "Itgb8"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3395593333333333
$ws.Range("H4").Value = 1.018678
$ws.Range("I4").Value = 0.6166145092460882
$ws.Range("J4").Value = 0.6166145092460882
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.673183333333334
$ws.Range("N4").Value = 26.01955
$ws.Range("O4").Value = 0.6868809929060228
$ws.Range("P4").Value = 0.6868809929060229
$ws.Range("Q4").Value = 2.945060350544445
$ws.Range("R4").Value = 26.5055431549
$ws.Range("S4").Value = 0.423540786351213
$ws.Range("T4").Value = 0.4235407863512131

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Lama1"
$ws.Range("C5").Value = "Itgb8"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.211124
$ws.Range("H5").Value = 0.633372
$ws.Range("I5").Value = 0.3833854907539118
$ws.Range("J5").Value = 0.3833854907539118
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1591403333333333
$ws.Range("N5").Value = 0.477421
$ws.Range("O5").Value = 0.01260326986877891
$ws.Range("P5").Value = 0.01260326986877891
$ws.Range("Q5").Value = 0.03359834373466666
$ws.Range("R5").Value = 0.302385093612
$ws.Range("S5").Value = 0.004831910803745791
$ws.Range("T5").Value = 0.004831910803745792

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Lama1"
$ws.Range("C6").Value = "Itgb8"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.211124
$ws.Range("H6").Value = 0.633372
$ws.Range("I6").Value = 0.3833854907539118
$ws.Range("J6").Value = 0.3833854907539118
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.794584666666667
$ws.Range("N6").Value = 11.383754
$ws.Range("O6").Value = 0.3005157372251983
$ws.Range("P6").Value = 0.3005157372251983
$ws.Range("Q6").Value = 0.8011278931653334
$ws.Range("R6").Value = 7.210151038488
$ws.Range("S6").Value = 0.1152133733953562
$ws.Range("T6").Value = 0.1152133733953563

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Lama1"
$ws.Range("C7").Value = "Itgb8"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.211124
$ws.Range("H7").Value = 0.633372
$ws.Range("I7").Value = 0.3833854907539118
$ws.Range("J7").Value = 0.3833854907539118
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.673183333333334
$ws.Range("N7").Value = 26.01955
$ws.Range("O7").Value = 0.6868809929060228
$ws.Range("P7").Value = 0.6868809929060229
$ws.Range("Q7").Value = 1.831117158066667
$ws.Range("R7").Value = 16.4800544226
$ws.Range("S7").Value = 0.2633402065548098
$ws.Range("T7").Value = 0.2633402065548098
